# Edit: Add "version" column (with value 1) everywhere, per commit message
# "Add version 1 everywhere": insert a new "version" column at the front of
# the "Export as TSV" sheet, backed by a new "version list" sheet (value "1"),
# mirroring the pattern used by the other *_list validation sheets.

function Add-ColLetter {
    param([string]$col, [int]$n)
    $num = 0
    foreach ($ch in $col.ToCharArray()) {
        $num = $num * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    $num += $n
    $result = ""
    while ($num -gt 0) {
        $rem = ($num - 1) % 26
        $result = [char]([int][char]'A' + $rem) + $result
        $num = [int](($num - 1) / 26)
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------------
# 1. Record existing data validations (sqref keyed by their CURRENT column)
#    so we can delete + recreate them after the insert, in the exact order
#    the target file uses (new "version" rule first).
# ---------------------------------------------------------------------------
$validationRules = @(
  @{col='I'; type=3; formula1="'assay_category list'!`$A`$1:`$A`$1"; title='Value must come from list'; msg='Value must be one of: sequence.'},
  @{col='J'; type=3; formula1="'assay_type list'!`$A`$1:`$A`$5"; title='Value must come from list'; msg='Value must be one of: scRNAseq-10xGenomics / scRNAseq / sciRNAseq / snRNAseq / SNARE2-RNAseq.'},
  @{col='K'; type=3; formula1="'analyte_class list'!`$A`$1:`$A`$1"; title='Value must come from list'; msg='Value must be one of: RNA.'},
  @{col='L'; type=3; formula1='"TRUE,FALSE"'; title='Not a boolean'; msg='The values in this column must be "TRUE" or "FALSE".'},
  @{col='AA'; type=3; formula1='"TRUE,FALSE"'; title='Not a boolean'; msg='The values in this column must be "TRUE" or "FALSE".'},
  @{col='AG'; type=2; formula1='-1e+307'; formula2='1e+307'; title='Not a number'; msg='The values in this column must be numbers.'},
  @{col='AH'; type=3; formula1="'library_final_yield_unit list'!`$A`$1:`$A`$1"; title='Value must come from list'; msg='Value must be one of: ng.'},
  @{col='AL'; type=2; formula1='-1e+307'; formula2='1e+307'; title='Not a number'; msg='The values in this column must be numbers.'},
  @{col='AM'; type=2; formula1='-1e+307'; formula2='1e+307'; title='Not a number'; msg='The values in this column must be numbers.'}
)

foreach ($r in $validationRules) {
    $ws.Range("$($r.col)2:$($r.col)1048576").Validation.Delete()
}

# ---------------------------------------------------------------------------
# 2. Record existing header-row comments (ref -> text) so we can move them to
#    the right-shifted column after the insert (comments do not auto-shift).
# ---------------------------------------------------------------------------
$oldComments = @(
    @{ref='A1'; text='HuBMAP Display ID of the donor of the assayed tissue.'},
    @{ref='B1'; text='HuBMAP Display ID of the assayed tissue.'},
    @{ref='C1'; text='Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.'},
    @{ref='D1'; text='DOI for protocols.io referring to the protocol for this assay.'},
    @{ref='E1'; text='Name of the person responsible for executing the assay.'},
    @{ref='F1'; text='Email address for the operator.'},
    @{ref='G1'; text='Name of the principal investigator responsible for the data.'},
    @{ref='H1'; text='Email address for the principal investigator.'},
    @{ref='I1'; text='Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.'},
    @{ref='J1'; text='The specific type of assay being executed.'},
    @{ref='K1'; text='Analytes are the target molecules being measured with the assay.'},
    @{ref='L1'; text='Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.'},
    @{ref='M1'; text='An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.'},
    @{ref='N1'; text='Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.'},
    @{ref='O1'; text='Link to a protocols document answering the question: How were single cells separated into a single-cell suspension?'},
    @{ref='P1'; text='The type of single cell entity derived from isolation protocol'},
    @{ref='Q1'; text='The method by which tissues are dissociated into single cells in suspension.'},
    @{ref='R1'; text='The method by which specific cell populations are sorted or enriched.'},
    @{ref='S1'; text='A quality metric by visual inspection prior to cell lysis or defined by known parameters such as wells with several cells or no cells. This can be captured at a high level.'},
    @{ref='T1'; text='Total number of cell/nuclei yielded post dissociation and enrichment'},
    @{ref='U1'; text='Number of cell/nuclei input to the assay'},
    @{ref='V1'; text='The kit used for the RNA sequencing assay'},
    @{ref='W1'; text='A link to the protocol document containing the library construction method (including version) that was used, e.g. "Smart-Seq2", "Drop-Seq", "10X v3".'},
    @{ref='X1'; text='Whether the library was generated for single-end or paired end sequencing'},
    @{ref='Y1'; text='Adapter sequence to be used for adapter trimming'},
    @{ref='Z1'; text='An id for the library. The id may be text and/or numbers'},
    @{ref='AA1'; text='Is the sequencing reaction run in repliucate, TRUE or FALSE'},
    @{ref='AB1'; text='Which read file contains the cell barcode'},
    @{ref='AC1'; text='Position(s) in the read at which the cell barcode starts.'},
    @{ref='AD1'; text='Length of the cell barcode in base pairs'},
    @{ref='AE1'; text='Number of PCR cycles to amplify cDNA'},
    @{ref='AF1'; text='Number of PCR cycles performed for library indexing'},
    @{ref='AG1'; text='Total number of ng of library after final pcr amplification step. This is the concentration (ng/ul) * volume (ul)'},
    @{ref='AH1'; text='Units of final library yield'},
    @{ref='AI1'; text='Average size of sequencing library fragments estimated via gel electrophoresis or bioanalyzer/tapestation.'},
    @{ref='AJ1'; text='Reagent kit used for sequencing'},
    @{ref='AK1'; text='Slash-delimited list of the number of sequencing cycles for, for example, Read1, i7 index, i5 index, and Read2.'},
    @{ref='AL1'; text='Percent of bases with Quality scores above Q30'},
    @{ref='AM1'; text='Percent PhiX loaded to the run'},
    @{ref='AN1'; text='Relative path to file with ORCID IDs for contributors for this dataset.'},
    @{ref='AO1'; text='Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.'}
)

foreach ($oc in $oldComments) {
    $ws.Range($oc.ref).Comment.Delete()
}

# ---------------------------------------------------------------------------
# 3. Insert the new "version" column at A, shifting every other column right
#    by one (cell values + data validations on untouched columns shift
#    automatically as part of Insert()).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# Header cell + style (copy header style from the neighboring cell so we
# reuse the existing bold/centered/wrap-text style instead of creating a new
# style entry).
$ws.Range("A1").Value = "version"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Re-create the header comments, shifted one column to the right, plus the
#    brand-new comment describing the "version" column.
# ---------------------------------------------------------------------------
$ws.Range("A1").AddComment("Current version of metadata schema. Template provides the correct value.")

foreach ($oc in $oldComments) {
    if ($oc.ref -match '^([A-Z]+)([0-9]+)$') {
        $colPart = $Matches[1]
        $rowPart = $Matches[2]
        $newCol = Add-ColLetter $colPart 1
        $newRef = "$newCol$rowPart"
        $ws.Range($newRef).AddComment($oc.text)
    }
}

# ---------------------------------------------------------------------------
# 5. Re-create the data validations on their shifted columns, plus the new
#    "version" validation first (matching target ordering).
# ---------------------------------------------------------------------------
$vRng = $ws.Range("A2:A1048576")
$vRng.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1")
$vRng.Validation.ErrorTitle = "Value must come from list"
$vRng.Validation.ErrorMessage = "Value must be one of: 1."

foreach ($r in $validationRules) {
    $newCol = Add-ColLetter $r.col 1
    $rng = $ws.Range("$($newCol)2:$($newCol)1048576")
    if ($r.type -eq 2) {
        $rng.Validation.Add(2, 1, 1, $r.formula1, $r.formula2)
    } else {
        $rng.Validation.Add(3, 1, 1, $r.formula1)
    }
    $rng.Validation.ErrorTitle = $r.title
    $rng.Validation.ErrorMessage = $r.msg
}

# ---------------------------------------------------------------------------
# 6. Add the new "version list" sheet right after "Export as TSV", holding
#    the single allowed value "1" (stored as text, matching the other *_list
#    sheets which store their allowed values as shared-string text).
# ---------------------------------------------------------------------------
$versionSheet = $wb.Worksheets.Add($null, $ws)
$versionSheet.Name = "version list"
$vCell = $versionSheet.Range("A1")
$vCell.NumberFormat = "@"
$vCell.Value = "1"
$vCell.Style = "Normal"

Write-Output "version column + version list sheet added"
